$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) — rows 3..6, column F ("想去人数")
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 191
$wsExpo.Range("F4").Value = 2395
$wsExpo.Range("F5").Value = 29
$wsExpo.Range("F6").Value = 526

# Sheet "全部类型" (all categories) — same events, offset by 2 rows — rows 5..8, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 191
$wsAll.Range("F6").Value = 2395
$wsAll.Range("F7").Value = 29
$wsAll.Range("F8").Value = 526
